$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.466.07'
$ws.Range('E2').Value = '  +0.69%  '
$ws.Range('D3').Value = '3.323.42'
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''586.13'
$ws.Range('E5').Value = '  +2.16%  '
$ws.Range('D6').Value = '''180.96'
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').Value = '''0.654'
$ws.Range('E7').Value = '  +6.02%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '3.319.72'
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('E10').Value = '  -0.35%  '
$ws.Range('D11').Value = '''6.83'
$ws.Range('E11').Value = '  +2.83%  '
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('D13').Value = '3.898.38'
$ws.Range('E13').Value = '  +0.22%  '
$ws.Range('E14').Value = '  -2.79%  '
$ws.Range('D15').Value = '66.488.65'
$ws.Range('E15').Value = '  +0.52%  '
$ws.Range('D16').Value = '''26.57'
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('E17').Value = '  -0.80%  '
$ws.Range('D18').Value = '3.297.20'
$ws.Range('E18').Value = '  -0.70%  '
$ws.Range('D19').Value = '''424.73'
$ws.Range('E19').Value = '  -2.50%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Value = '''13.17'
$ws.Range('E20').Value = '  -2.50%  '
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').Value = '''5.50'
$ws.Range('E21').Value = '  -2.88%  '
$ws.Range('E22').Value = '  -2.34%  '
$ws.Range('D23').Value = '''71.84'
$ws.Range('E23').Value = '  -1.89%  '
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('D26').Value = '3.471.08'
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('D27').Value = '''0.515'
$ws.Range('E27').Value = '  -0.69%  '
$ws.Range('E28').Value = '  +6.79%  '
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('D30').Value = '''9.09'
$ws.Range('E30').Value = '  +0.61%  '
$ws.Range('D31').Value = '''0.999'
$ws.Range('E31').Value = '  +0.57%  '
$ws.Range('E32').Value = '  -1.33%  '
$ws.Range('D33').Value = '''22.38'
$ws.Range('E33').Value = '  -1.24%  '
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('E35').Value = '  -0.79%  '
$ws.Range('D36').Value = '''6.62'
$ws.Range('E36').Value = '  -1.74%  '
$ws.Range('D37').Value = '''1.19'
$ws.Range('E37').Value = '  -1.82%  '
$ws.Range('D38').Value = '''160.18'
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('E39').Value = '  -2.31%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '''1.81'
$ws.Range('E40').Value = '  +1.25%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '2.867.81'
$ws.Range('E41').Value = '  +1.33%  '
$ws.Range('D42').Value = '''26.40'
$ws.Range('E42').Value = '  -4.37%  '
$ws.Range('E43').Value = '  -1.89%  '
$ws.Range('E44').Value = '  -3.62%  '
$ws.Range('D45').Value = '''39.80'
$ws.Range('E45').Value = '  -0.95%  '
$ws.Range('E46').Value = '  -0.22%  '
$ws.Range('D47').Value = '''5.91'
$ws.Range('E47').Value = '  -4.06%  '
$ws.Range('D48').Value = '''2.32'
$ws.Range('E48').Value = '  -0.58%  '
$ws.Range('D49').Value = '''23.20'
$ws.Range('E49').Value = '  -3.39%  '
$ws.Range('D50').Value = '''314.04'
$ws.Range('E50').Value = '  -2.34%  '
$ws.Range('E51').Value = '  +1.05%  '
